$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell formatting (bold, centered, bordered) from H1
# onto the two new header cells so they match the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-32 (columns I and J)
$data = @{
    2  = @(9, 9)
    3  = @(7, 8)
    4  = @(5, 7)
    5  = @(9, 9)
    6  = @(8, 9)
    7  = @(8, 8)
    8  = @(6, 6)
    9  = @(4, 6)
    10 = @(6, 8)
    11 = @(8, 9)
    12 = @(8, 9)
    13 = @(8, 9)
    14 = @(9, 9)
    15 = @(5, 6)
    16 = @(9, 9)
    17 = @(7, 9)
    18 = @(8, 9)
    19 = @(1, 5)
    20 = @(1, 1)
    21 = @(1, 5)
    22 = @(1, 6)
    23 = @(1, 6)
    24 = @(1, 5)
    25 = @(1, 6)
    26 = @(1, 5)
    27 = @(1, 7)
    28 = @(1, 4)
    29 = @(3, 7)
    30 = @(5, 7)
    31 = @(5, 6)
    32 = @(1, 2)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
